$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = "2º Teste"
$ws.Range("C4").Value = "novo commit"

$ws.Range("C5").Select()
